$p = $ppt.ActivePresentation

# Remove the first slide (the REST-verb diagram slide); the remaining
# slide becomes slide 1.
$p.Slides.Item(1).Delete()

$s = $p.Slides.Item(1)

# Re-apply the "GET" label text on its textbox. Re-typing it through
# Delete+InsertAfter (instead of a plain TextRange.Text assignment)
# mirrors the author re-touching that label, which drops the stale,
# redundant trailing paragraph-mark run properties left over from a
# previous edit of this textbox.
$shp = $s.Shapes.Item(6)
$shp.TextFrame.TextRange.Delete()
$shp.TextFrame.TextRange.InsertAfter("GET")
